$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("U3").Value = 44692
$ws.Range("V3").Value = 44693
Write-Output "done"
